# Rename worksheets and update values per diff

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename sheets
$ws1.Name = "data_RAM code_FLASH"
$ws2.Name = "data_RAM code_CCM"

# Sheet 1 (data_RAM code_FLASH) - intensity row (row 2)
$ws1.Range("B2").Value = 13105
$ws1.Range("C2").Value = 23163
$ws1.Range("D2").Value = 31151

# Sheet 1 (data_RAM code_FLASH) - energy row (row 5)
$ws1.Range("B5").Value = 6.142
$ws1.Range("C5").Value = 6.334
$ws1.Range("D5").Value = 6.848

# Sheet 2 (data_RAM code_CCM) - intensity row (row 2)
$ws2.Range("B2").Value = 11093
$ws2.Range("C2").Value = 21438
$ws2.Range("D2").Value = 31586

# Sheet 2 (data_RAM code_CCM) - energy row (row 5)
$ws2.Range("B5").Value = 5.199
$ws2.Range("C5").Value = 5.036
$ws2.Range("D5").Value = 4.971
